$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" column header at H1, copying the formatting of the
# neighboring header cell (G1) so it keeps the same bold/border/centered
# style used by the rest of the header row.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Values for the new "Save" column (rows 2-12)
$saveValues = @(0, 0, 1, 0, 1, 1, 0, 0, 1, 1, 0)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
